# Added New Mac-Address and Document Types
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch workbook calculation to manual (calcPr calcMode="manual")
$excel.Calculation = -4135

# Append new data row 33 (mirrors the pattern of the preceding rows)
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"
$ws.Range("H33").Value = "now()"

# Update the view: scroll back to top-left default and move the selection to E31
[void]$ws.Range("A1").Select()
[void]$ws.Range("E31").Select()
